$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rename the "Then_Question" / "Else_Question" header cells to
# "Then_Goto" / "Else_Goto" respectively.
$ws.Range("I1").Value = "Then_Goto"
$ws.Range("J1").Value = "Else_Goto"

# Match the author's resulting selection (active cell moved to I1).
$ws.Range("I1").Select()
